$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the two claim-number cells used for the smart folders
# (values are stored as text with a quote-prefix, so force text entry)
$ws.Range("F2").Value = "'0420194406833"
$ws.Range("F3").Value = "'0420172008637  "

# Move the current selection to F6
$ws.Range("F6").Select()
